# Deploying to gh-pages — add the 2022 column (S) to the transport stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S; since S is immediately past the existing table
# (which currently ends at R), this duplicates column R's cell formatting
# into the new column S for every row that has data in R, without disturbing
# any existing columns.
$ws.Columns("S").Insert()

# Year header.
$ws.Range("S3").Value = 2022

# Data rows (2022 figures), row by row.
$ws.Range("S4").Value = 10444.200000000001
$ws.Range("S5").Value = 21.7
$ws.Range("S6").Value = 7361.6
$ws.Range("S7").Value = 143.1
$ws.Range("S8").Value = 844.2
# Row 9 has no 2022 figure; leave the cell blank (already created by the
# column insert above with the correct formatting).
$ws.Range("S10").Value = "2 756,0"
$ws.Range("S11").Value = "1 013,8"
$ws.Range("S12").Value = "1 451,1"
$ws.Range("S13").Value = 273.39999999999998
$ws.Range("S14").Value = "-"
$ws.Range("S15").Value = 17.7

# Match the author's recorded selection at save time.
$ws.Range("T3").Select()
